$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while keeping it as literal text (matches the
# source workbook, which stores every cell - including numeric-looking
# prices - as inline strings) and leaving the cell style untouched.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Row 2: 'Bitcoin'
Set-TextValue $ws.Range('D2') '64.273.28'
$ws.Range('E2').Value = '  +1.08%  '

# Row 3: 'Ethereum'
Set-TextValue $ws.Range('D3') '3.501.99'
$ws.Range('E3').Value = '  +0.58%  '

# Row 4: 'TetherUSD'
$ws.Range('E4').Value = '  -0.20%  '

# Row 5: 'BNB'
Set-TextValue $ws.Range('D5') '586.46'
$ws.Range('E5').Value = '  +0.82%  '

# Row 6: 'Solana'
Set-TextValue $ws.Range('D6') '134.22'
$ws.Range('E6').Value = '  +2.52%  '

# Row 7: 'USDC'
$ws.Range('E7').Value = '  -0.04%  '

# Row 8: 'XRP'
$ws.Range('E8').Value = '  +0.36%  '

# Row 9: 'Dogecoin'
$ws.Range('E9').Value = '  +2.61%  '

# Row 10: 'Toncoin'
Set-TextValue $ws.Range('D10') '7.27'
$ws.Range('E10').Value = '  +2.02%  '

# Row 11: 'Cardano'
$ws.Range('E11').Value = '  +1.58%  '

# Row 12: 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range('D12') '4.096.01'
$ws.Range('E12').Value = '  -0.08%  '

# Row 13: 'TRON'
$ws.Range('E13').Value = '  +1.07%  '

# Row 14: 'ShibaInu'
$ws.Range('E14').Value = '  +3.70%  '

# Row 15: 'WrappedEther'
Set-TextValue $ws.Range('D15') '3.500.86'
$ws.Range('E15').Value = '  -0.15%  '

# Row 16: 'Avalanche'
Set-TextValue $ws.Range('D16') '26.14'
$ws.Range('E16').Value = '  -4.27%  '

# Row 17: 'WrappedBTC'
Set-TextValue $ws.Range('D17') '64.368.45'
$ws.Range('E17').Value = '  +0.87%  '

# Row 18: 'Uniswap'
Set-TextValue $ws.Range('D18') '9.92'
$ws.Range('E18').Value = '  +0.15%  '

# Row 19: 'Polkadot'
Set-TextValue $ws.Range('D19') '5.75'
$ws.Range('E19').Value = '  +1.99%  '

# Row 20: 'Chainlink'
Set-TextValue $ws.Range('D20') '13.81'
$ws.Range('E20').Value = '  -3.15%  '

# Row 21: 'BitcoinCash'
Set-TextValue $ws.Range('D21') '392.26'
$ws.Range('E21').Value = '  +2.33%  '

# Row 22: 'Polygon'
$ws.Range('E22').Value = '  -0.49%  '

# Row 23: 'WrappedeETH'
Set-TextValue $ws.Range('D23') '3.641.29'
$ws.Range('E23').Value = '  +0.40%  '

# Row 24: 'Litecoin'
Set-TextValue $ws.Range('D24') '74.20'
$ws.Range('E24').Value = '  +2.05%  '

# Row 25: 'Dai'
$ws.Range('E25').Value = '  +0.07%  '

# Row 26: 'PEPE'
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range('D26') '0.0000115'
$ws.Range('E26').Value = '  +2.44%  '

# Row 27: 'RenderToken'
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D27') '7.57'
$ws.Range('E27').Value = '  +1.69%  '

# Row 28: 'Binance-PegBSC-USD'
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range('D28') '0.997'
$ws.Range('E28').Value = '  -0.31%  '

# Row 29: 'Fetch.AI'
$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D29') '1.52'
$ws.Range('E29').Value = '  -3.55%  '

# Row 30: 'InternetComputer(DFINITY)'
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D30') '8.29'
$ws.Range('E30').Value = '  +0.69%  '

# Row 31: 'PancakeSwap'
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D31') '2.23'
$ws.Range('E31').Value = '  -0.67%  '

# Row 32: 'RenzoRestakedETH'
$ws.Range('B32').Value = 'RenzoRestakedETH'
$ws.Range('C32').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextValue $ws.Range('D32') '3.523.69'
$ws.Range('E32').Value = '  +0.95%  '

# Row 33: 'USDe'
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range('D33') '1.00'
$ws.Range('E33').Value = '  +0.00%  '

# Row 34: 'Kaspa'
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D34') '0.149'
$ws.Range('E34').Value = '  +3.64%  '

# Row 35: 'EthereumClassic'
$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D35') '23.48'
$ws.Range('E35').Value = '  +0.21%  '

# Row 36: 'NEARProtocol'
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D36') '5.21'
$ws.Range('E36').Value = '  -1.51%  '

# Row 37: 'ImmutableX'
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D37') '1.57'
$ws.Range('E37').Value = '  +0.34%  '

# Row 38: 'Aptos'
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D38') '6.92'
$ws.Range('E38').Value = '  +0.05%  '

# Row 39: 'Monero'
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D39') '161.87'
$ws.Range('E39').Value = '  -0.20%  '

# Row 40: 'Hedera'
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D40') '0.0783'
$ws.Range('E40').Value = '  -1.28%  '

# Row 41: 'Mantle'
$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D41') '0.807'
$ws.Range('E41').Value = '  -0.06%  '

# Row 42: 'EnergySwap'
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D42') '25.56'
$ws.Range('E42').Value = '  -2.85%  '

# Row 43: 'FirstDigitalUSD'
$ws.Range('E43').Value = '  -0.34%  '

# Row 44: 'Filecoin'
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D44') '4.41'
$ws.Range('E44').Value = '  +1.03%  '

# Row 45: 'ONDO'
$ws.Range('B45').Value = 'ONDO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue $ws.Range('D45') '1.18'
$ws.Range('E45').Value = '  -2.95%  '

# Row 46: 'Stacks'
$ws.Range('E46').Value = '  +1.89%  '

# Row 47: 'Maker'
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D47') '2.469.15'
$ws.Range('E47').Value = '  +2.43%  '

# Row 48: 'Cosmos'
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D48') '6.80'
$ws.Range('E48').Value = '  +0.01%  '

# Row 49: 'SuiNetwork'
$ws.Range('B49').Value = 'SuiNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range('D49') '0.897'
$ws.Range('E49').Value = '  +1.09%  '

# Row 50: 'VeChain'
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D50') '0.0262'
$ws.Range('E50').Value = '  -0.62%  '

# Row 51: 'TheGraph'
$ws.Range('B51').Value = 'TheGraph'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range('D51') '0.217'
$ws.Range('E51').Value = '  -0.93%  '
